$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.802.29'
$ws.Range("E2").Value = '  -1.21%  '

$ws.Range("D3").Value = '1.942.47'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.15'
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4894'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2944'
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06909'
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.41'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '106.34'
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = '1.956.63'
$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07728'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.347'
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6999'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '276.97'
$ws.Range("E16").Value = '  -3.32%  '
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = '30.797.30'
$ws.Range("E17").Value = '  -0.86%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007726'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '13.12'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D20").Style = "Normal"

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.191.93'
$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.443'
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.528'
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.710'
$ws.Range("E25").Value = '  -2.85%  '
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.34'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.66'
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.156'
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1042'
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = '  -4.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.553'
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.547'
$ws.Range("E32").Value = '  -5.64%  '
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = '  -3.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04855'
$ws.Range("E34").Value = '  -3.38%  '
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7496'
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.157'
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9997'
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.724'
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01986'
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.661'
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.61'
$ws.Range("E41").Value = '  +6.85%  '
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.447'
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.092'
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9069'
$ws.Range("E44").Value = '  +2.24%  '
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = '  -1.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4403'
$ws.Range("E46").Value = '  -1.27%  '
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9979'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.733'
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '985.37'
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1243'
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.275'
$ws.Range("E51").Value = '  -1.22%  '
$ws.Range("D51").Style = "Normal"
